# Daily attendance processing - reorder "Recorded By" (column G) tokens.
# For every data row, the comma-separated list of recorders in column G
# is reversed in order (e.g. "a, b" -> "b, a").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2

    if ($val -ne $null -and $val -ne "") {
        $parts = $val -split ", "
        $n = $parts.Count
        if ($n -gt 1) {
            $reversedParts = @()
            for ($i = $n - 1; $i -ge 0; $i--) {
                $reversedParts += $parts[$i]
            }
            $newVal = $reversedParts -join ", "
            $cell.Value2 = $newVal
        }
    }
}
